$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 47500
$ws.Range("I10").Value = 47500
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 47500
$ws.Range("L10").ClearContents()
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("H64").Value = 22824.25
$ws.Range("I64").Value = 30306
$ws.Range("J64").Value = 7860.75
$ws.Range("K64").Value = 30306
$ws.Range("L64").Value = 7860.75
$ws.Range("M64").Value = -30058
$ws.Range("N64").Value = -8356.75
$ws.Range("H67").Value = 22824.25
$ws.Range("I67").Value = 30306
$ws.Range("J67").Value = 7860.75
$ws.Range("K67").Value = 30306
$ws.Range("L67").Value = 7860.75
$ws.Range("M67").Value = -29448
$ws.Range("N67").Value = -9576.75
$ws.Range("H74").Value = 7342.5
$ws.Range("I74").Value = 6963
$ws.Range("J74").Value = 9999
$ws.Range("K74").Value = 6963
$ws.Range("L74").Value = 9999
$ws.Range("M74").Value = -6027
$ws.Range("N74").Value = -11871
$ws.Range("H77").Value = 7342.5
$ws.Range("I77").Value = 6963
$ws.Range("J77").Value = 9999
$ws.Range("K77").Value = 34815
$ws.Range("L77").Value = 49995
$ws.Range("M77").Value = -30135
$ws.Range("N77").Value = -59355
$ws.Range("H87").Value = 150000
$ws.Range("J87").Value = 150000
$ws.Range("L87").Value = 150000
$ws.Range("N87").Value = -152496
$ws.Range("H90").Value = 150000
$ws.Range("J90").Value = 150000
$ws.Range("L90").Value = 450000
$ws.Range("N90").Value = -462480
$ws.Range("H101").Value = 3195
$ws.Range("J101").Value = 3195
$ws.Range("L101").Value = 9585
$ws.Range("N101").Value = -12829
$ws.Range("H113").Value = 6743.923
$ws.Range("I113").Value = 5937.1
$ws.Range("J113").Value = 9433.333000000001
$ws.Range("K113").Value = 5937.1
$ws.Range("L113").Value = 9433.333000000001
$ws.Range("M113").Value = -2683.1
$ws.Range("N113").Value = -15941.333
$ws.Range("H132").Value = 22090.059
$ws.Range("I132").Value = 25362.242
$ws.Range("K132").Value = 76086.726
$ws.Range("M132").Value = -73556.726
$ws.Range("H137").Value = 36231.832
$ws.Range("I137").Value = 24533.715
$ws.Range("J137").Value = 77175.25
$ws.Range("K137").Value = 73601.145
$ws.Range("L137").Value = 231525.75
$ws.Range("M137").Value = -71051.145
$ws.Range("N137").Value = -236625.75
$ws.Range("H138").Value = 49913.043
$ws.Range("I138").Value = 4692.385
$ws.Range("J138").Value = 108699.9
$ws.Range("K138").Value = 14077.155
$ws.Range("L138").Value = 326099.7
$ws.Range("M138").Value = -8937.155000000001
$ws.Range("N138").Value = -336379.7
$ws.Range("H141").Value = 1340.0625
$ws.Range("I141").Value = 1340.0625
$ws.Range("K141").Value = 4020.1875
$ws.Range("M141").Value = 1159.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1142.7097
$ws.Range("I2").Value = 1179.7916
$ws.Range("J2").Value = 1015.5714
$ws.Range("K2").Value = 1179.7916
$ws.Range("L2").Value = 1015.5714
$ws.Range("M2").Value = -1066.7916
$ws.Range("N2").Value = -1241.5714
$ws.Range("H32").Value = 14822.013
$ws.Range("I32").Value = 15030.27
$ws.Range("J32").Value = 6700
$ws.Range("K32").Value = 15030.27
$ws.Range("L32").Value = 6700
$ws.Range("M32").Value = -14743.27
$ws.Range("N32").Value = -7274
$ws.Range("H74").Value = 134492.56
$ws.Range("I74").Value = 182955.7
$ws.Range("J74").Value = 11470.77
$ws.Range("K74").Value = 182955.7
$ws.Range("L74").Value = 11470.77
$ws.Range("M74").Value = -182081.7
$ws.Range("N74").Value = -13218.77
$ws.Range("H77").Value = 134492.56
$ws.Range("I77").Value = 182955.7
$ws.Range("J77").Value = 11470.77
$ws.Range("K77").Value = 914778.5
$ws.Range("L77").Value = 57353.85000000001
$ws.Range("M77").Value = -910410.5
$ws.Range("N77").Value = -66089.85000000001
$ws.Range("H116").Value = 1142.7097
$ws.Range("I116").Value = 1179.7916
$ws.Range("J116").Value = 1015.5714
$ws.Range("K116").Value = 1179.7916
$ws.Range("L116").Value = 1015.5714
$ws.Range("M116").Value = 1114.2084
$ws.Range("N116").Value = -5603.5714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1142.7097
$ws.Range("I3").Value = 1179.7916
$ws.Range("J3").Value = 1015.5714
$ws.Range("K3").Value = 1179.7916
$ws.Range("L3").Value = 1015.5714
$ws.Range("M3").Value = -1065.7916
$ws.Range("N3").Value = -1243.5714
$ws.Range("H22").Value = 512.5714
$ws.Range("I22").Value = 512.5714
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 512.5714
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H86").Value = 749.25
$ws.Range("I86").Value = 749.25
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 749.25
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 749.25
$ws.Range("I89").Value = 749.25
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 3746.25
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("H14").Value = 5000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("H31").Value = 1697178
$ws.Range("I31").Value = 2382955
$ws.Range("J31").Value = 2905.353
$ws.Range("K31").Value = 2382955
$ws.Range("L31").Value = 2905.353
$ws.Range("M31").Value = -2382660
$ws.Range("N31").Value = -3495.353
$ws.Range("H34").Value = 1697178
$ws.Range("I34").Value = 2382955
$ws.Range("J34").Value = 2905.353
$ws.Range("K34").Value = 2382955
$ws.Range("L34").Value = 2905.353
$ws.Range("M34").Value = -2382753
$ws.Range("N34").Value = -3309.353
$ws.Range("H132").Value = 51077.75
$ws.Range("I132").Value = 63146.25
$ws.Range("J132").Value = 2803.75
$ws.Range("K132").Value = 189438.75
$ws.Range("L132").Value = 8411.25
$ws.Range("M132").Value = -186908.75
$ws.Range("N132").Value = -13471.25
$ws.Range("H134").Value = 1898.129
$ws.Range("I134").Value = 1640.7858
$ws.Range("K134").Value = 4922.357400000001
$ws.Range("M134").Value = -2387.357400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4490.0586
$ws.Range("J68").Value = 4614.4375
$ws.Range("L68").Value = 13843.3125
$ws.Range("N68").Value = -15465.3125
$ws.Range("H71").Value = 4490.0586
$ws.Range("J71").Value = 4614.4375
$ws.Range("L71").Value = 41529.9375
$ws.Range("N71").Value = -49641.9375
$ws.Range("H131").Value = 2321.923
$ws.Range("I131").Value = 3505.1538
$ws.Range("J131").Value = 1730.3077
$ws.Range("K131").Value = 10515.4614
$ws.Range("L131").Value = 5190.9231
$ws.Range("M131").Value = -5475.4614
$ws.Range("N131").Value = -15270.9231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 6000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 6000
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("H97").Value = 1262.0333
$ws.Range("I97").Value = 1209.5555
$ws.Range("J97").Value = 1340.75
$ws.Range("K97").Value = 1209.5555
$ws.Range("L97").Value = 1340.75
$ws.Range("M97").Value = -713.5554999999999
$ws.Range("N97").Value = -2332.75
$ws.Range("H132").Value = 1816.2075
$ws.Range("I132").Value = 1572.5476
$ws.Range("K132").Value = 4717.642800000001
$ws.Range("M132").Value = -2187.642800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2268.25
$ws.Range("I55").Value = 1613.6666
$ws.Range("J55").Value = 3359.2222
$ws.Range("K55").Value = 1613.6666
$ws.Range("L55").Value = 3359.2222
$ws.Range("M55").Value = -1440.6666
$ws.Range("N55").Value = -3705.2222
$ws.Range("H136").Value = 4433
$ws.Range("I136").Value = 4432.5386
$ws.Range("J136").Value = 4435
$ws.Range("K136").Value = 13297.6158
$ws.Range("L136").Value = 13305
$ws.Range("M136").Value = -10747.6158
$ws.Range("N136").Value = -18405

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1922.5
$ws.Range("I62").Value = 1896.6666
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 1896.6666
$ws.Range("L62").Value = 2000
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 1922.5
$ws.Range("I65").Value = 1896.6666
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 9483.333000000001
$ws.Range("L65").Value = 10000
$ws.Range("N65").Value = -16240
$ws.Range("H96").Value = 76924344
$ws.Range("I96").Value = 1375
$ws.Range("J96").Value = 1000000000
$ws.Range("K96").Value = 1375
$ws.Range("L96").Value = 1000000000
$ws.Range("M96").Value = -2
$ws.Range("N96").Value = -1000002746
$ws.Range("H126").Value = 158144.9
$ws.Range("I126").Value = 1616.2222
$ws.Range("K126").Value = 4848.6666
$ws.Range("M126").Value = -2378.6666
$ws.Range("H132").Value = 42053.777
$ws.Range("I132").Value = 49604.266
$ws.Range("J132").Value = 4301.3335
$ws.Range("K132").Value = 148812.798
$ws.Range("L132").Value = 12904.0005
$ws.Range("M132").Value = -146282.798
$ws.Range("N132").Value = -17964.0005
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
